$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in row 2 (E2:T2 partial) to the new TPM-derived values
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 0.026348808596
$ws.Range("R2").Value = 0.237139277364
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove row 3 (the Resolving-Mac -> Alb/Lrp2/MuSCs entry) entirely
$ws.Rows("3:3").Delete()
